$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 2772.04168703704
$ws.Range("C2").Value = 509.651133176672
$ws.Range("D2").Value = 4593.572127043365
